$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("total_sonuc_keywords")
$ws.Activate()
$ws.Range("F:G").EntireColumn.Delete()
